# "New data: length of stay" -- append a new continuous-outcome block
# (LOS = Length of Stay) to the bottom of the extraction sheet, mirroring
# the layout already used by the other outcome blocks (header/divider row,
# then one row per study with Dex/Control mean+-SD, sample sizes and a
# formatted "mean ± sd" label pair).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Section header / divider row (row 65) ---------------------------
# Same study-name label as the very first divider row of the sheet
# (shared string "Shi et al"), spanning B:J with a top border.
$ws.Range("B65:J65").Value = "Shi et al"
$ws.Range("B65:J65").Borders.Item(8).LineStyle = 1
$ws.Range("B65:J65").Borders.Item(8).Weight = 2

# Row 66 is a spacer row that just carries the next study label.
$ws.Range("B66").Value = "Xie et al"

# ---- Data rows (67-72): LOS mean +- SD per study ----------------------
$ws.Range("B67").Value = "Meng et al"
$ws.Range("C67").Value = 4.1
$ws.Range("D67").Value = 1.48
$ws.Range("E67").Value = 4.6
$ws.Range("F67").Value = 1.88
$ws.Range("G67").Value = 20
$ws.Range("H67").Value = 20
$ws.Range("I67").Value = "LOS"
$ws.Range("J67").Value = "4.10 ± 1.48"
$ws.Range("K67").Value = "4.60 ± 1.88"

$ws.Range("B68").Value = "Zhang et al"
$ws.Range("C68").Value = 7.1
$ws.Range("D68").Value = 3.9
$ws.Range("E68").Value = 7.76
$ws.Range("F68").Value = 3.12
$ws.Range("G68").Value = 28
$ws.Range("H68").Value = 28
$ws.Range("I68").Value = "LOS"
$ws.Range("J68").Value = "7.10 ± 3.90"
$ws.Range("K68").Value = "7.76 ± 3.12"

$ws.Range("B69").Value = "Wu et al"
$ws.Range("C69").Value = 5.6
$ws.Range("D69").Value = 2.5
$ws.Range("E69").Value = 5.9
$ws.Range("F69").Value = 2.5
$ws.Range("G69").Value = 30
$ws.Range("H69").Value = 30
$ws.Range("I69").Value = "LOS"
$ws.Range("J69").Value = "5.60 ± 2.50"
$ws.Range("K69").Value = "5.90 ± 2.50"

$ws.Range("B70").Value = "Jannu et al"
$ws.Range("C70").Value = 6.33
$ws.Range("D70").Value = 4.61
$ws.Range("E70").Value = 9.66
$ws.Range("F70").Value = 5.38
$ws.Range("G70").Value = 40
$ws.Range("H70").Value = 40
$ws.Range("I70").Value = "LOS"
$ws.Range("J70").Value = "6.33 ± 4.61"
$ws.Range("K70").Value = "9.66 ± 5.38"

$ws.Range("B71").Value = "Lee et al (1)"
$ws.Range("C71").Value = 7.16
$ws.Range("D71").Value = 3.93
$ws.Range("E71").Value = 8.03
$ws.Range("F71").Value = 3.93
$ws.Range("G71").Value = 25
$ws.Range("H71").Value = 25
$ws.Range("I71").Value = "LOS"
$ws.Range("J71").Value = "7.16 ± 3.93"
$ws.Range("K71").Value = "8.03 ± 3.93"

$ws.Range("B72").Value = "Lee et al (2)"
$ws.Range("C72").Value = 6.1
$ws.Range("D72").Value = 4.58
$ws.Range("E72").Value = 7.13
$ws.Range("F72").Value = 3.81
$ws.Range("G72").Value = 50
$ws.Range("H72").Value = 50
$ws.Range("I72").Value = "LOS"
$ws.Range("J72").Value = "6.10 ± 4.58"
$ws.Range("K72").Value = "7.13 ± 3.81"

# Trailing spacer row with the last study label, closing out the block.
$ws.Range("B73").Value = "Zhou et al"

# ---- View state: scroll the new block into view and select it --------
$ws.Range("J65:K73").Select()
